$d = $word.ActiveDocument

# 1. Date in the header line: 15.06.24 -> 14.06.24
$d.Paragraphs(1).Range.Text = "⚡️🚀המאמר היומי של מייק 14.06.24:⚡️🚀"

# 2. Paper title
$d.Paragraphs(2).Range.Text = "CLLMs: Consistency Large Language Models"

# 3. Intro paragraph (replaced with new content about iterative parallel decoding methods)
$d.Paragraphs(3).Range.Text = "בשתי הסקירות הקודמות(כדאי שתעברו עליהם כי נתתי שם קצת הסברים) דיברנו על שיטות איטרטיביות מקבילות לדגימה ממודלי שפה. השיטות האלו מבוססות על שיטות יאקובי או (Gauss-Seidel (GS. השיטות האלו מתחילות מכמות מסוימת n של טוקנים שנדגמים באקראי (או בצורה קצת יותר מושכלת) ואז מעדכנים טוקנים אלו בבת אחת באיטרציות עד שתנאי עצירה מתקיים(התכנסות). תנאי העצירה כאן הוא בד״כ שוויון בין הפלטים של איטרציות עוקבות. "

# 4. Next paragraph
$d.Paragraphs(4).Range.Text = "מובן שאנו מעוניינים לסיים את התהליך במשמעות פחות איטרציות ממספר הטוקנים שאנו חוזים בו זמנית (ד״א ניתן להראות נדרשות לכל היותר ח איטרציות עד ההתכנסות). "

# 5. Next paragraph
$d.Paragraphs(5).Range.Text = "שימו לב שמהלך האימון של מודלי שפה מותאם לשיטת הדגימה האוטו-רגרסיביות כאשר בוחרים טוקן בעל הסתסברות הגבוה ביותר ביהנתן הטוקנים הקודמים. אולם עכשיו אנו דוגמים בצורה אחרת ואולי ניתן להתחשב בזה במהלך האימון. כלומר במהלך האימון אשכרה דוגמים עם השיטה הזו (השילוב של יאקובי ו- GS)."

# 6. Next paragraph
$d.Paragraphs(6).Range.Text = "וזה בדיוק מה שנסקור אותו היום עושה. המחברים מוסיפים עוד איבר ללוס הרגיל של מודלי שפה (הממקסם את הנראות המירבית של הדאטה). מטרת האיבר הזה היא לגרום למזעור של מספר האיטרציות עד להתכנסות של הדגימה האיטרטיבית. "

# 7. Next paragraph
$d.Paragraphs(7).Range.Text = "המחברים בחנו שתי אופציות לאיבר הזה:"

# Insert three new paragraphs after paragraph 7 (before the URL paragraph)
$p7 = $d.Paragraphs(7)
$p7.Range.InsertParagraphAfter()
$d.Paragraphs(8).Range.Text = "מזעור של מרחק (KL הפוך לדעתי אך לא צללתי לעומק) בין התפלגויות הטוקנים בנקודת ההתכנסות לבין התפלגויות טוקנים במהלך הדגימה האיטרטיבית (דוגמים האיטרציות באקראי)."

$p8 = $d.Paragraphs(8)
$p8.Range.InsertParagraphAfter()
$d.Paragraphs(9).Range.Text = "מזעור מרחק בין התפלגויות הטוקנים באיטרציות עוקבות."

$p9 = $d.Paragraphs(9)
$p9.Range.InsertParagraphAfter()
$d.Paragraphs(10).Range.Text = "ואם חשבתם שיש דמיון בין השיטה הזו לבין המאמר של איליה סלוצקבר ושותפיו ""Consistency Models"" - אכן הוא קיים ואני אצלול בו בקרוב."

# 11. Final paragraph: replace the old arXiv link with the new one
$d.Paragraphs(11).Range.Text = "https://arxiv.org/abs/2403.00835"
